$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - copy style from H1 (bold/centered/bordered header)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# New data cells
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9
$ws.Range("I3").Value = 10
$ws.Range("J3").Value = 10
